# Apply crypto price/volume updates per commit "Updated cryptos list" (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Value)
    # Force the written value to be stored as literal text, exactly as scraped,
    # instead of letting Excel auto-parse numeric-looking strings into numbers.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "27.056.03"
Set-TextValue $ws.Cells.Item(2, 5) "  -2.36%  "

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "1.821.22"
Set-TextValue $ws.Cells.Item(3, 5) "  -1.51%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 5) "  -1.40%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "311.02"
Set-TextValue $ws.Cells.Item(5, 5) "  -2.68%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "1.0000"
Set-TextValue $ws.Cells.Item(6, 5) "  -1.20%  "

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "0.4220"
Set-TextValue $ws.Cells.Item(7, 5) "  -2.11%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "0.3679"
Set-TextValue $ws.Cells.Item(8, 5) "  -1.94%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 5) "  -1.87%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "0.8417"
Set-TextValue $ws.Cells.Item(10, 5) "  -4.13%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 5) "  -3.92%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "1.816.53"
Set-TextValue $ws.Cells.Item(12, 5) "  -1.79%  "

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "6.656"
Set-TextValue $ws.Cells.Item(13, 5) "  -1.23%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "0.07056"
Set-TextValue $ws.Cells.Item(14, 5) "  -1.22%  "

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "5.280"
Set-TextValue $ws.Cells.Item(15, 5) "  -2.95%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "89.80"
Set-TextValue $ws.Cells.Item(16, 5) "  +0.53%  "

# Row 17
Set-TextValue $ws.Cells.Item(17, 5) "  -1.45%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "0.000008784"
Set-TextValue $ws.Cells.Item(18, 5) "  -2.43%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "1.000"
Set-TextValue $ws.Cells.Item(19, 5) "  -1.14%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 5) "  -3.86%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "27.120.81"
Set-TextValue $ws.Cells.Item(21, 5) "  -2.19%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "5.126"
Set-TextValue $ws.Cells.Item(22, 5) "  -1.71%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "10.85"
Set-TextValue $ws.Cells.Item(23, 5) "  -2.30%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "2.046.39"
Set-TextValue $ws.Cells.Item(24, 5) "  -1.63%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "1.979"
Set-TextValue $ws.Cells.Item(25, 5) "  -1.22%  "

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "151.82"
Set-TextValue $ws.Cells.Item(26, 5) "  -2.30%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 5) "  +2.30%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "18.25"
Set-TextValue $ws.Cells.Item(28, 5) "  -2.22%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) "5.274"
Set-TextValue $ws.Cells.Item(29, 5) "  -1.95%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "116.27"
Set-TextValue $ws.Cells.Item(30, 5) "  -2.52%  "

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) "0.08745"
Set-TextValue $ws.Cells.Item(31, 5) "  -2.27%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "1.178"
Set-TextValue $ws.Cells.Item(32, 5) "  -4.19%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) "0.7396"
Set-TextValue $ws.Cells.Item(33, 5) "  -5.01%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) "2.947"
Set-TextValue $ws.Cells.Item(34, 5) "  +0.64%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "4.417"
Set-TextValue $ws.Cells.Item(35, 5) "  -3.24%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "0.9997"
Set-TextValue $ws.Cells.Item(36, 5) "  -1.36%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "1.092"
Set-TextValue $ws.Cells.Item(37, 5) "  -3.69%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "0.01949"
Set-TextValue $ws.Cells.Item(38, 5) "  -1.59%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "0.05251"
Set-TextValue $ws.Cells.Item(39, 5) "  -1.99%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "7.356"
Set-TextValue $ws.Cells.Item(40, 5) "  +0.05%  "

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "2.872"
Set-TextValue $ws.Cells.Item(41, 5) "  -0.77%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "0.1691"
Set-TextValue $ws.Cells.Item(42, 5) "  -0.33%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "0.5041"
Set-TextValue $ws.Cells.Item(43, 5) "  -2.11%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "8.594"
Set-TextValue $ws.Cells.Item(44, 5) "  -2.68%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "10.47"
Set-TextValue $ws.Cells.Item(45, 5) "  -2.21%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "106.32"
Set-TextValue $ws.Cells.Item(46, 5) "  -2.02%  "

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) "0.4715"
Set-TextValue $ws.Cells.Item(47, 5) "  -1.78%  "

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) "0.9994"
Set-TextValue $ws.Cells.Item(48, 5) "  -1.34%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "RenderToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Cells.Item(49, 4) "1.898"
Set-TextValue $ws.Cells.Item(49, 5) "  +2.44%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Cells.Item(50, 4) "0.06338"
Set-TextValue $ws.Cells.Item(50, 5) "  -2.27%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "1.649"
Set-TextValue $ws.Cells.Item(51, 5) "  -2.70%  "
